$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new coding-scheme row entries (abbreviation + description pair)
$ws.Range("A13").Value = "Reorg"
$ws.Range("B13").Value = "Re-organized data"

# Match the formatting of the rest of that row (C13:D13 use cell style index 3)
$ws.Range("C13:D13").Copy()
$ws.Range("A13:B13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the active selection left by the author after finishing the edit
$ws.Range("B11").Select()
